$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.442.05'
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.807.80'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.83'
$ws.Range("E5").Value = '  -0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.603'

$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.32'
$ws.Range("E8").Value = '  +3.86%  '

$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0681'
$ws.Range("E10").Value = '  -1.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0965'
$ws.Range("E11").Value = '  +1.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.062.84'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.33'
$ws.Range("E13").Value = '  +1.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.818.56'
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.631'
$ws.Range("E15").Value = '  -1.53%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.45'
$ws.Range("E16").Value = '  +2.73%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.401.54'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.66'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.95'
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0{0}0774" -f [char]0x2083
$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  -2.04%  '

$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("E24").Value = '  +5.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.16'
$ws.Range("E25").Value = '  +0.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.88'
$ws.Range("E26").Value = '  +4.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.36'
$ws.Range("E27").Value = '  +3.97%  '

$ws.Range("E28").Value = '  +2.13%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.80'
$ws.Range("E30").Value = '  -0.36%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.92'
$ws.Range("E31").Value = '  -1.42%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0517'
$ws.Range("E33").Value = '  -2.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("E34").Value = '  -1.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.362.11'
$ws.Range("E35").Value = '  -2.36%  '

$ws.Range("E36").Value = '  -4.11%  '

$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -6.22%  '

$ws.Range("E39").Value = '  -2.03%  '

$ws.Range("E40").Value = '  +1.64%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '81.08'
$ws.Range("E41").Value = '  -2.16%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  -1.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.940'
$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("E44").Value = '  +5.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.39'
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("E46").Value = '  -2.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.966.69'
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.80'
$ws.Range("E48").Value = '  -3.07%  '

$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.56'
$ws.Range("E50").Value = '  -1.89%  '

$ws.Range("E51").Value = '  -6.36%  '
